$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B8").Value = 5.715999999999999
$ws.Range("B10").Value = 6.944
$ws.Range("B12").Value = 6.444
$ws.Range("C13").Value = -12.729
$ws.Range("B18").Value = 6.873
$ws.Range("B25").Value = 6.629
